# Update LR-pairs sheet with new TPM-derived numbers and a new "Inflammatory-Mac"
# sending-cluster row (inserted before the existing MuSCs row, which moves down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Make room: insert a new row at position 3 so the former row 3 (MuSCs ->
#    MuSCs) becomes row 4, and row 3 is free for the new Inflammatory-Mac row.
# ---------------------------------------------------------------------------
$ws.Rows.Item(3).Insert()

# ---------------------------------------------------------------------------
# 2) Refresh the text/label columns (A-D) for rows 2-4 so the shared-string
#    table is rebuilt with "Inflammatory-Mac" placed right after "FAPs".
#    Clear first, then re-enter values in the order: FAPs, Inflammatory-Mac,
#    MuSCs, Artn, Gfra3.
# ---------------------------------------------------------------------------
$ws.Range("A2:D4").Value2 = ""

$ws.Cells.Item(2,1).Value2 = "FAPs"
$ws.Cells.Item(3,1).Value2 = "Inflammatory-Mac"
$ws.Cells.Item(4,1).Value2 = "MuSCs"

$ws.Cells.Item(2,4).Value2 = "MuSCs"
$ws.Cells.Item(3,4).Value2 = "MuSCs"
$ws.Cells.Item(4,4).Value2 = "MuSCs"

$ws.Cells.Item(2,2).Value2 = "Artn"
$ws.Cells.Item(3,2).Value2 = "Artn"
$ws.Cells.Item(4,2).Value2 = "Artn"

$ws.Cells.Item(2,3).Value2 = "Gfra3"
$ws.Cells.Item(3,3).Value2 = "Gfra3"
$ws.Cells.Item(4,3).Value2 = "Gfra3"

# ---------------------------------------------------------------------------
# 3) Row 2 (FAPs -> MuSCs): keep E-G/K-L/O-P the same, update recalculated
#    values.
# ---------------------------------------------------------------------------
$ws.Cells.Item(2,5).Value2 = 3
$ws.Cells.Item(2,6).Value2 = 1
$ws.Cells.Item(2,7).Value2 = 1.415594
$ws.Cells.Item(2,8).Value2 = 4.246782
$ws.Cells.Item(2,9).Value2 = 0.7538443241324221
$ws.Cells.Item(2,10).Value2 = 0.7538443241324221
$ws.Cells.Item(2,11).Value2 = 2
$ws.Cells.Item(2,12).Value2 = 0.6666666666666666
$ws.Cells.Item(2,13).Value2 = 0.186728
$ws.Cells.Item(2,14).Value2 = 0.560184
$ws.Cells.Item(2,15).Value2 = 1
$ws.Cells.Item(2,16).Value2 = 1
$ws.Cells.Item(2,17).Value2 = 0.264331036432
$ws.Cells.Item(2,18).Value2 = 2.378979327888
$ws.Cells.Item(2,19).Value2 = 0.7538443241324221
$ws.Cells.Item(2,20).Value2 = 0.7538443241324221

# ---------------------------------------------------------------------------
# 4) Row 3 (new Inflammatory-Mac -> MuSCs) values.
# ---------------------------------------------------------------------------
$ws.Cells.Item(3,5).Value2 = 1
$ws.Cells.Item(3,6).Value2 = 0.3333333333333333
$ws.Cells.Item(3,7).Value2 = 0.2356576666666667
$ws.Cells.Item(3,8).Value2 = 0.706973
$ws.Cells.Item(3,9).Value2 = 0.1254944528268394
$ws.Cells.Item(3,10).Value2 = 0.1254944528268394
$ws.Cells.Item(3,11).Value2 = 2
$ws.Cells.Item(3,12).Value2 = 0.6666666666666666
$ws.Cells.Item(3,13).Value2 = 0.186728
$ws.Cells.Item(3,14).Value2 = 0.560184
$ws.Cells.Item(3,15).Value2 = 1
$ws.Cells.Item(3,16).Value2 = 1
$ws.Cells.Item(3,17).Value2 = 0.04400388478133333
$ws.Cells.Item(3,18).Value2 = 0.396034963032
$ws.Cells.Item(3,19).Value2 = 0.1254944528268394
$ws.Cells.Item(3,20).Value2 = 0.1254944528268394

# ---------------------------------------------------------------------------
# 5) Row 4 (former MuSCs -> MuSCs row) values.
# ---------------------------------------------------------------------------
$ws.Cells.Item(4,5).Value2 = 2
$ws.Cells.Item(4,6).Value2 = 0.6666666666666666
$ws.Cells.Item(4,7).Value2 = 0.2265816666666667
$ws.Cells.Item(4,8).Value2 = 0.679745
$ws.Cells.Item(4,9).Value2 = 0.1206612230407385
$ws.Cells.Item(4,10).Value2 = 0.1206612230407385
$ws.Cells.Item(4,11).Value2 = 2
$ws.Cells.Item(4,12).Value2 = 0.6666666666666666
$ws.Cells.Item(4,13).Value2 = 0.186728
$ws.Cells.Item(4,14).Value2 = 0.560184
$ws.Cells.Item(4,15).Value2 = 1
$ws.Cells.Item(4,16).Value2 = 1
$ws.Cells.Item(4,17).Value2 = 0.04230914145333334
$ws.Cells.Item(4,18).Value2 = 0.38078227308
$ws.Cells.Item(4,19).Value2 = 0.1206612230407385
$ws.Cells.Item(4,20).Value2 = 0.1206612230407385

Write-Host "Edit complete"
